$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.164924666666667
$ws.Range("H2").Value = 3.494774
$ws.Range("M2").Value = 0.2887656666666666
$ws.Range("N2").Value = 0.866297
$ws.Range("O2").Value = 0.02693956104254192
$ws.Range("P2").Value = 0.02693956104254193
$ws.Range("Q2").Value = 0.3363902479864444
$ws.Range("R2").Value = 3.027512231878
$ws.Range("S2").Value = 0.02693956104254192
$ws.Range("T2").Value = 0.02693956104254193

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.164924666666667
$ws.Range("H3").Value = 3.494774
$ws.Range("O3").Value = 0.687775905612258
$ws.Range("P3").Value = 0.6877759056122581
$ws.Range("Q3").Value = 8.588154316347332
$ws.Range("R3").Value = 77.293388847126
$ws.Range("S3").Value = 0.687775905612258
$ws.Range("T3").Value = 0.6877759056122581

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.164924666666667
$ws.Range("H4").Value = 3.494774
$ws.Range("M4").Value = 3.057970333333333
$ws.Range("N4").Value = 9.173911
$ws.Range("O4").Value = 0.2852845333452001
$ws.Range("P4").Value = 0.2852845333452002
$ws.Range("Q4").Value = 3.562305071234889
$ws.Range("R4").Value = 32.060745641114
$ws.Range("S4").Value = 0.2852845333452001
$ws.Range("T4").Value = 0.2852845333452002
